# RTS-96.xlsx — "Infeasibilities in bigger scenarios"
# The "Costo" (cost, column H) values on the Parametros sheet were stored
# in raw currency units; rescale them down by 1e6 (e.g. 38388000 -> 38.388)
# to fix infeasibilities that show up in larger scenarios.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parametros")

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 8)   # column H = 8
    $cell.Value = $cell.Value2 / 1000000
}
